$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value = "r775"
$ws.Cells.Item($row, 2).Value = "fred"
$ws.Cells.Item($row, 3).Value = "4:16 fred"
$ws.Cells.Item($row, 4).Value = "2025-10-01 16:16:21"
